$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.569.00"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "3.034.67"
$ws.Range("E3").Value = "  +2.50%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'384.26"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").Value = "'102.59"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").Value = "'0.544"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.586"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").Value = "'36.80"
$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "'0.0859"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").Value = "3.517.57"
$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("D14").Value = "'18.66"
$ws.Range("E14").Value = "  +1.79%  "

$ws.Range("D15").Value = "'7.74"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").Value = "3.042.92"
$ws.Range("E16").Value = "  +2.43%  "

$ws.Range("D17").Value = "'0.975"
$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("D18").Value = "'10.61"
$ws.Range("E18").Value = "  -10.89%  "

$ws.Range("D19").Value = "51.603.29"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").Value = "'3.08"
$ws.Range("E20").Value = "  -0.45%  "

$ws.Range("D21").Value = "'12.35"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'69.91"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").Value = "'266.85"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("E25").Value = "  -3.52%  "

$ws.Range("D26").Value = "'8.34"
$ws.Range("E26").Value = "  +5.15%  "

$ws.Range("D27").Value = "'7.47"
$ws.Range("E27").Value = "  +4.16%  "

$ws.Range("E28").Value = "  +4.63%  "

$ws.Range("D29").Value = "'26.31"
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "'0.107"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").Value = "'10.28"
$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("D34").Value = "'34.06"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("D35").Value = "'50.51"
$ws.Range("E35").Value = "  -1.07%  "

$ws.Range("D36").Value = "'0.0447"
$ws.Range("E36").Value = "  +2.58%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "'3.38"
$ws.Range("E38").Value = "  +3.71%  "

$ws.Range("D39").Value = "'0.285"
$ws.Range("E39").Value = "  +6.06%  "

$ws.Range("D40").Value = "'16.99"
$ws.Range("E40").Value = "  +2.40%  "

$ws.Range("D41").Value = "'1.86"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("D42").Value = "'0.116"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("D43").Value = "'127.23"
$ws.Range("E43").Value = "  +2.14%  "

$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("D45").Value = "'3.67"
$ws.Range("E45").Value = "  +3.87%  "

$ws.Range("D46").Value = "'21.62"
$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").Value = "'2.48"
$ws.Range("E47").Value = "  +3.01%  "

$ws.Range("D48").Value = "'2.10"
$ws.Range("E48").Value = "  +3.68%  "

$ws.Range("D49").Value = "2.040.20"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").Value = "3.339.61"
$ws.Range("E50").Value = "  +2.80%  "

$ws.Range("E51").Value = "  +6.07%  "
